# Rename the severity-level headers (E1:L1) on every sheet and drop the
# now-redundant "Category" (upper_primary) column M that followed them.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("E1").Value() = "% severity levels 1-2"
    $ws.Range("F1").Value() = "# severity levels 1-2"
    $ws.Range("G1").Value() = "% severity level 3"
    $ws.Range("H1").Value() = "# severity level 3"
    $ws.Range("I1").Value() = "% severity level 4"
    $ws.Range("J1").Value() = "# severity level 4"
    $ws.Range("K1").Value() = "% severity level 5"
    $ws.Range("L1").Value() = "# severity level 5"

    $ws.Range("M:M").EntireColumn.Delete()
}
